# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for each localized-language sheet once a handback has been
# received, and flips the synced rows' status from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere that status is shown.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR for RGB 6495ED - matches the workbook's existing "HyperLink" cell style
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/790a6c8ed4cf5afb6faca0ac8f6b928fdcebbcec/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet: both language columns move from "Ready for handoff" to
# "Handed back: in sync with en-US" for each of the two tracked files.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): mark status handed back, record the
# generated handback package + target file + handback timestamp.
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; K2 = "2016-08-25 00:27:27"; K3 = "2016-08-25 00:27:27" },
    @{ Name = "de-de"; K2 = "2016-08-25 00:27:34"; K3 = "2016-08-25 00:27:34" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) for both source-file rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Row 2 -> 2550924d-c725-4b1b-83d4-7dc6865b6c7f.md
    $ws.Hyperlinks.Add($ws.Range("I2"), ($ghBase + "2550924d-c725-4b1b-83d4-7dc6865b6c7f.md"), "", "", "2550924d-c725-4b1b-83d4-7dc6865b6c7f.md")
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("J2").Value = ("2550924d-c725-4b1b-83d4-7dc6865b6c7f.e0e46d739bb8247967148d1fdaa6db185f2cd904." + $lang.Name + ".xlf")
    $ws.Range("K2").Value = $lang.K2

    # Row 3 -> 4b06aa87-bf77-4c8b-8cd2-802b4ee32c8b.md
    $ws.Hyperlinks.Add($ws.Range("I3"), ($ghBase + "4b06aa87-bf77-4c8b-8cd2-802b4ee32c8b.md"), "", "", "4b06aa87-bf77-4c8b-8cd2-802b4ee32c8b.md")
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = $hyperlinkColor
    $ws.Range("J3").Value = ("4b06aa87-bf77-4c8b-8cd2-802b4ee32c8b.bfcb4ea5c7c0446f26bce03b1aaac43b4c35e08e." + $lang.Name + ".xlf")
    $ws.Range("K3").Value = $lang.K3

    $ws.Columns.Item(3).ColumnWidth = 29.2
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
